$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 413 entirely ("「毎年、ノーベル賞が授けられる」" post) so that
# all subsequent rows shift up by one (row 414 becomes 413, ..., row 587
# becomes 586), matching the author's edit.
$ws.Rows.Item(413).Delete()
